# Append three new log rows (255-257) to Sheet1, matching the diff:
#   A255/B255 = 2023-12-12 20:27:49 / 0.002
#   A256/B256 = 2023-12-12 20:28:40 / 0.0016
#   A257/B257 = 2023-12-12 20:30:10 / 0.0048

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2023-12-12 20:27:49", 0.002),
    @("2023-12-12 20:28:40", 0.0016),
    @("2023-12-12 20:30:10", 0.0048)
)

$startRow = 255
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
